# Update crypto "Volume(1h)" percentage values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '  -2.37%  '
$ws.Range("E3").Value = '  -5.70%  '
$ws.Range("E4").Value = '  +1.09%  '
$ws.Range("E5").Value = '  -1.69%  '
$ws.Range("E6").Value = '  -8.08%  '
$ws.Range("E7").Value = '  -0.80%  '
$ws.Range("E8").Value = '  -11.33%  '
$ws.Range("E9").Value = '  -5.13%  '
$ws.Range("E10").Value = '  -3.02%  '
$ws.Range("E11").Value = '  +0.03%  '
$ws.Range("E12").Value = '  -2.89%  '
$ws.Range("E13").Value = '  -3.59%  '
$ws.Range("E14").Value = '  -7.16%  '
$ws.Range("E15").Value = '  -5.29%  '
$ws.Range("E16").Value = '  -6.21%  '
$ws.Range("E17").Value = '  -20.36%  '
$ws.Range("E18").Value = '  -9.73%  '
$ws.Range("E19").Value = '  -4.12%  '
$ws.Range("E20").Value = '  -2.70%  '
$ws.Range("E21").Value = '  -2.46%  '
$ws.Range("E22").Value = '  -7.65%  '
$ws.Range("E23").Value = '  -0.36%  '
$ws.Range("E24").Value = '  -2.24%  '
$ws.Range("E25").Value = '  -10.36%  '
$ws.Range("E26").Value = '  +5.58%  '
$ws.Range("E27").Value = '  -0.22%  '
$ws.Range("E28").Value = '  -5.36%  '
$ws.Range("E29").Value = '  -5.47%  '
$ws.Range("E30").Value = '  -7.45%  '
$ws.Range("E31").Value = '  -8.75%  '
$ws.Range("E32").Value = '  -12.41%  '
$ws.Range("E33").Value = '  -2.41%  '
$ws.Range("E34").Value = '  -6.32%  '
$ws.Range("E35").Value = '  -9.01%  '
$ws.Range("E36").Value = '  -0.38%  '
$ws.Range("E37").Value = '  -5.13%  '
$ws.Range("E38").Value = '  +0.86%  '
$ws.Range("E39").Value = '  -3.17%  '
$ws.Range("E40").Value = '  -10.20%  '
$ws.Range("E41").Value = '  -1.58%  '
$ws.Range("E42").Value = '  +0.34%  '
$ws.Range("E43").Value = '  -0.10%  '
$ws.Range("E44").Value = '  -1.59%  '
$ws.Range("E45").Value = '  -3.30%  '
$ws.Range("E46").Value = '  -2.99%  '
$ws.Range("E47").Value = '  -11.50%  '
$ws.Range("E48").Value = '  -4.34%  '
$ws.Range("E49").Value = '  -11.93%  '
$ws.Range("E50").Value = '  -4.58%  '
$ws.Range("E51").Value = '  -4.87%  '
